$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 30-32, column A: set value "yes" and apply style "4" (red font / green fill)
# so they match the formatting already used by the "yes" cells elsewhere (e.g. A2),
# while B:E in these rows switch from style 3 to style 4 as well.
$ws.Range("A30:E32").Style = "Style4"

$ws.Range("A30").Value = "yes"
$ws.Range("A31").Value = "yes"
$ws.Range("A32").Value = "yes"

# Update the sheet view: scroll back to the top (remove frozen/topLeft offset) and
# change the active selection cell.
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("E18").Select()
